$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update flanking-gene (nearest upstream/downstream ORF) values for a few EVEs
$ws.Range("G15").Value = "ENSTGUG00000000478"
$ws.Range("H15").Value = "LRRC2"
$ws.Range("G17").Value = "ENSTGUG00000004184 (intronic)"
$ws.Range("H17").Value = "ENSTGUG00000004184 (intronic)"
$ws.Range("H18").Value = "MGST1"
$ws.Range("H20").Value = "LIG3"
$ws.Range("G21").Value = "ATP2B2 (intronic)"
$ws.Range("H21").Value = "ATP2B2 (intronic)"

# Column width tweaks (A, C, F, G got narrower)
$ws.Columns.Item(1).ColumnWidth = 27.666666666666668
$ws.Columns.Item(3).ColumnWidth = 17.330729166666668
$ws.Columns.Item(6).ColumnWidth = 8.166666666666666
$ws.Columns.Item(7).ColumnWidth = 26.830729166666668

# Move the active selection to where the user ended up working
$ws.Range("H21").Select()
